$wb = $excel.ActiveWorkbook

# Update the "Status" value from "Ready for handoff" to "In Translation"
# wherever it appears: Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The status column width shrinks to fit the new (shorter) text
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
